$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("B2").Value = "https://go.mwater.co/terre_neuve_cap"
$ws.Range("D2").Value = "Commune Action Plan: https://go.mwater.co/terre_neuve_cap"

$ws.Range("B3").Value = "---"
$ws.Range("D3").Value = "Project performance: ---"
